# Add season-record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from AC1) onto the three new header cells,
# then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill every data row (2-66) with the team's season record.
for ($row = 2; $row -le 66; $row++) {
    $ws.Cells.Item($row, 30).Value = 77
    $ws.Cells.Item($row, 31).Value = 85
    $ws.Cells.Item($row, 32).Value = 0
}
